$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Common")

# Step 1: remove all existing comments in the affected region (rows 86-127)
# so stale comments do not linger at the old anchor rows after the shift.
for ($r = 86; $r -le 127; $r++) {
    $cmt = $ws.Cells.Item($r, 1).Comment
    if ($cmt -ne $null) {
        $cmt.Delete()
    }
}

# Step 2: insert the two new category header rows.
# First insertion splits the old "KVM RAM" block into "KVM and VCenter RAM" (VSD/VSC/VSTAT)
# and a new "KVM RAM" block (VCIN/NUH/Webfilter/Portal).
$ws.Rows.Item(89).Insert()
# Second insertion (row 99, post first-shift numbering) splits the old "CPU" block into
# "KVM and VCenter CPU" (VSD/VSC/VSTAT/VNSUTIL) and a new "KVM CPU" block (NUH/VCIN/Portal/Webfilter).
$ws.Rows.Item(99).Insert()

# Step 3: give the two new header rows the same look as the other category headers
# (style copied from the existing "KVM RAM" header at row 85), then merge A:B and set text.
$ws.Range("A85:B85").Copy()
$ws.Range("A89:B89").PasteSpecial(-4122)
$ws.Range("A89").Value = "KVM RAM"
$ws.Range("A89:B89").Merge()

$ws.Range("A85:B85").Copy()
$ws.Range("A99:B99").PasteSpecial(-4122)
$ws.Range("A99").Value = "KVM CPU"
$ws.Range("A99:B99").Merge()
$excel.CutCopyMode = 0

# Step 4: rename cells whose label text changed.
$ws.Range("A85").Value = 'KVM and VCenter RAM'
$ws.Range("A86").Value = 'KVM VSD RAM'
$ws.Range("A87").Value = 'KVM VSC RAM'
$ws.Range("A88").Value = 'KVM VSTAT RAM'
$ws.Range("A94").Value = 'KVM and VCenter CPU'
$ws.Range("A95").Value = 'KVM VSD CPU cores'
$ws.Range("A96").Value = 'KVM VSC CPU cores'
$ws.Range("A97").Value = 'KVM VSTAT CPU cores'
$ws.Range("A98").Value = 'KVM VNSUTIL CPU cores'

# Step 5: re-create every comment at its correct final row with its final text
# (text is unchanged for most; seven of them get the reworded "For KVM and VCenter
# deployments: ..." phrasing per the commit).
$comments = @(
    @{Row=86; Text='For KVM and VCenter deployments: amount of VSD RAM to allocate, in GB. Note: Values smaller than the default are for lab and PoC only. Production deployments must use a value greater than or equal to the default. [default: 24]'},
    @{Row=87; Text='For KVM and VCenter deployments: amount of VSC RAM to allocate, in GB. Note: Values smaller than the default are for lab and PoC only. Production deployments must use a value greater than or equal to the default. [default: 4]'},
    @{Row=88; Text='For KVM and VCenter deployments: amount of VSTAT RAM to allocate, in GB. Note: Values smaller than the default are for lab and PoC only. Production deployments must use a value greater than or equal to the default. [default: 16]'},
    @{Row=90; Text='Amount of VCIN RAM to allocate, in GB. Valid only for KVM deployments. Note: Values smaller than the default are for lab and PoC only. Production deployments must use a value greater than or equal to the default. [default: 24]'},
    @{Row=91; Text='Amount of NUH RAM to allocate, in GB. Valid only for KVM deployments. Note: Values smaller than the default are for lab and PoC only. Production deployments must use a value greater than or equal to the default. [default: 8]'},
    @{Row=92; Text='Amount of Webfilter RAM to allocate, in GB. Valid only for KVM deployments. Note: Values smaller than the default are for lab and PoC only. Production deployments must use a value greater than or equal to the default. [default: 8]'},
    @{Row=93; Text='Amount of Portal RAM to allocate, in GB. Valid only for KVM deployments. Note: Values smaller than the default are for lab and PoC only. Production deployments must use a value greater than or equal to the default. [default: 24]'},
    @{Row=95; Text='For KVM and VCenter deployments: number of CPU''s for VSD. [default: 6]'},
    @{Row=96; Text='For KVM and VCenter deployments: number of CPU''s for VSC. [default: 6]'},
    @{Row=97; Text='For KVM and VCenter deployments: number of CPU''s for VSTAT. [default: 6]'},
    @{Row=98; Text='For KVM and VCenter deployments: number of CPU''s for VNSUTIL. [default: 2]'},
    @{Row=100; Text='Number of CPU''s for NUH. Valid only for KVM deployments [default: 2]'},
    @{Row=101; Text='Number of CPU''s for VCIN. Valid only for KVM deployments [default: 6]'},
    @{Row=102; Text='Number of CPU''s for Portal vm. Valid only for KVM deployments [default: 6]'},
    @{Row=103; Text='Number of CPU''s for Webfilter vm. Valid only for KVM deployments [default: 2]'},
    @{Row=105; Text='VSD Architect URL. Required for tasks during Upgrade, Health Checks etc [default: https://(vsd_fqdn):8443]'},
    @{Row=106; Text='Enterprise name used for authentication with VSD Architect. Required for tasks during Upgrade, Health Checks etc [default: csp]'},
    @{Row=107; Text='VCIN URL used for API interaction. Required for tasks like VRS-E upgrade (through VCIN) [default: https://(vcin_ip_address):8443]'},
    @{Row=108; Text='Enterprise name used for authentication with VCIN. Required for tasks like VRS-E upgrade (through VCIN) [default: csp]'},
    @{Row=110; Text='List of hooks files (List items separated by comma.)'},
    @{Row=111; Text='Skip tasks and playbooks (List items separated by comma.)'},
    @{Row=113; Text='Address of SMTP server to be used if emailing health results'},
    @{Row=114; Text='Port to be used on the SMTP Server [default: 25]'},
    @{Row=115; Text='Email address from which health report will be sent'},
    @{Row=116; Text='List of destination email addresses (List items separated by comma.)'},
    @{Row=118; Text='Address of the mail server to be used to receive monit alerts via email'},
    @{Row=119; Text='Port on mail server to be used for monit alerts [default: 25]'},
    @{Row=120; Text='Encryption to be used when sending monit alerts via email'},
    @{Row=121; Text='Enables use of monit eventqueue to store alerts if email alerts fail to send [default: True]'},
    @{Row=122; Text='Email address from which monit alerts will be sent'},
    @{Row=123; Text='Email address to reply to monit alert emails'},
    @{Row=124; Text='Email subject for alert emails. Overrides monit default alert subject'},
    @{Row=125; Text='Email message for alert emails. Overrides monit default alert message'},
    @{Row=126; Text='Destination email address for monit alerts'},
    @{Row=127; Text='Specific events for which alerts should be sent. One string can be used to hold multiple events, separated by commas'},
    @{Row=128; Text='Events for which alerts should not be sent. One string can be used to hold multiple events, separated by commas'},
    @{Row=129; Text='Allowing VSD in-place upgrade during Installation [default: False]'}
)
foreach ($item in $comments) {
    $cell = $ws.Cells.Item($item.Row, 1)
    [void]$cell.AddComment($item.Text)
}

